# Auto-generated edit script: updates Price (D) and Volume(1h) (E) columns
# with refreshed cryptos data, matching the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.123.59"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").Value = "1.632.23"
$ws.Range("E3").Value = "  -0.88%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.93%  "
$ws.Range("E6").Value = "  +1.09%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  -1.31%  "
$ws.Range("E9").Value = "  -0.69%  "
$ws.Range("E10").Value = "  -0.58%  "
$ws.Range("E11").Value = "  +0.16%  "
$ws.Range("D12").Value = "1.860.10"
$ws.Range("E12").Value = "  -0.93%  "
$ws.Range("D13").Value = "1.636.57"
$ws.Range("E13").Value = "  -1.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.12"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.540"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.62%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.64"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.13%  "
$ws.Range("D17").Value = "27.098.92"
$ws.Range("E17").Value = "  -0.29%  "
$ws.Range("D18").Value = "0.0₃0732"
$ws.Range("E18").Value = "  -1.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "214.26"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.83"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.19%  "
$ws.Range("E22").Value = "  -0.91%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.83%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.38"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("E28").Value = "  -1.14%  "
$ws.Range("E29").Value = "  -1.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0504"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.28%  "
$ws.Range("E31").Value = "  -0.69%  "
$ws.Range("E32").Value = "  +0.36%  "
$ws.Range("E33").Value = "  -0.95%  "
$ws.Range("D34").Value = "1.304.19"
$ws.Range("E34").Value = "  +2.19%  "
$ws.Range("E35").Value = "  -0.98%  "
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("E37").Value = "  -1.48%  "
$ws.Range("E38").Value = "  +0.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.843"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("E40").Value = "  -0.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.26"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.69%  "
$ws.Range("E42").Value = "  -0.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.27"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.87%  "
$ws.Range("D44").Value = "1.769.42"
$ws.Range("E44").Value = "  -1.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.26"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.98%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "90.63"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.86%  "
$ws.Range("E47").Value = "  +0.52%  "
$ws.Range("E48").Value = "  +0.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.815"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +20.97%  "
$ws.Range("E50").Value = "  -0.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.57"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.24%  "
